# Deploy the implementation guide.
# - Bump Status from "draft" to "active" and update the publication Date
#   on the Metadata sheet.
# - Add the ele-1/ext-1 invariant text to the root Extension row's
#   Invariants column (AJ1) on the Elements sheet - it was already present
#   on the Extension.extension row (AJ3) but missing from the root row.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B6").Value = "active"
$meta.Range("B8").Value = "2023-10-16T18:33:36+00:00"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
# Setting a multi-line value on row 1 makes the engine auto-grow the row
# height (wrapText is on for this style); AutoFit brings it back in line
# with the sheet default so no spurious <row ht=.../> attribute is emitted.
$elements.Rows.Item(1).AutoFit()
